$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 6): add headers for the new columns.
# Order matters for shared-string table indices: precision, fitness, log-based precision
$ws.Range("L6").Value = "precision"
$ws.Range("Q6").Value = "fitness"
$ws.Range("N6").Value = "log-based precision"

# Column L: "precision" raw values (rows 7-24, skipping the separator row 15)
$precision = @{
    7  = 0.375
    8  = 0.4318
    9  = 0.375
    10 = 0.3809
    11 = 0.3928
    12 = 0.3125
    13 = 0.4318
    14 = 0.3333
    16 = 0.4772
    17 = 0.333
    18 = 0
    19 = 0.2
    20 = 0.4772
    21 = 0.2692
    22 = 0.4772
    23 = 0.2
    24 = 0.2
}

# Column Q: "fitness" raw values (rows 7-24, skipping the separator row 15)
$fitness = @{
    7  = 0.99
    8  = 0.8968
    9  = 0.9905
    10 = 0.7825
    11 = 0.90037
    12 = 0.521
    13 = 0.8989
    14 = 1
    16 = 0.8744
    17 = 0.988
    18 = 0.93333
    19 = 0.987
    20 = 0.8751
    21 = 0.9471
    22 = 0.8756
    23 = 0.9294
    24 = 0.987
}

foreach ($row in (7..24 | Where-Object { $_ -ne 15 })) {
    $ws.Range("L$row").Value = $precision[$row]
    $ws.Range("Q$row").Value = $fitness[$row]
}

# Column N: "log-based precision" = precision * 2
# N7 is its own (non-shared) formula; N8:N24 form one shared-formula group,
# matching the original author's fill-down pattern.
$ws.Range("N7").Formula = "=L7*2"
$ws.Range("N8:N24").Formula = "=L8*2"

# Row 15 is a "NULL" separator row (no numeric data available for that log).
$ws.Range("L15").Value = "NULL"
$ws.Range("N15").Value = "NULL"

# Leave the selection where the author last clicked while reviewing the new data.
$null = $ws.Range("P15").Select()
